$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Remove the old "com.td.asignado" row (row 125). Everything below shifts
#    up by one (row 126 -> 125, row 127 -> 126, ... row 172 -> 171), which is
#    exactly what the target diff shows (row N now holds what used to be in
#    row N+1).
# ---------------------------------------------------------------------------
$ws.Rows("125").Delete()

# ---------------------------------------------------------------------------
# 2) Append the two new translation rows at the bottom of the table
#    (new rows 172 and 173), reusing the existing D/E formula pattern
#    (shared formula referencing $A, $B/$C and $D$1/$E$1).
# ---------------------------------------------------------------------------
$dFormulaTemplate = $ws.Range("D171").Formula
$eFormulaTemplate = $ws.Range("E171").Formula

$ws.Range("A172").Value = "com.td.seguro"
$ws.Range("B172").Value = "¿Está seguro?"
$ws.Range("C172").Value = "Are you sure?"
$ws.Range("D172").Formula = $dFormulaTemplate.Replace("171", "172")
$ws.Range("E172").Formula = $eFormulaTemplate.Replace("171", "172")

$ws.Range("A173").Value = "com.td.tutor.asignado"
$ws.Range("B173").Value = "El tutor esta asignado"
$ws.Range("C173").Value = "The tutor is assigned"
$ws.Range("D173").Formula = $dFormulaTemplate.Replace("171", "173")
$ws.Range("E173").Formula = $eFormulaTemplate.Replace("171", "173")

# ---------------------------------------------------------------------------
# 3) Update the sheet view to match where the author ended up: scrolled down
#    near the bottom of the table with the newly-typed row selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D172:E172").Select()
$excel.ActiveWindow.ScrollRow = 166
$excel.ActiveWindow.ScrollColumn = 4
